$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.445.81"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").Value = "3.475.05"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'554.54"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'179.63"
$ws.Range("E6").Value = "  -4.20%  "
$ws.Range("E7").Value = "  +3.63%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "'54.04"
$ws.Range("E11").Value = "  -4.91%  "
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "4.029.86"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "'18.74"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.475.98"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "65.427.19"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'416.85"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +3.61%  "
$ws.Range("D23").Value = "'85.97"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'12.84"
$ws.Range("E25").Value = "  +8.23%  "
$ws.Range("E26").Value = "  -9.44%  "
$ws.Range("D27").Value = "'2.86"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").Value = "'6.03"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "  +4.66%  "
$ws.Range("D30").Value = "'30.34"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -5.77%  "
$ws.Range("D32").Value = "'613.12"
$ws.Range("E32").Value = "  -10.62%  "
$ws.Range("D33").Value = "'11.78"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "'59.15"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'37.59"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.145"
$ws.Range("E38").Value = "  +8.27%  "
$ws.Range("D39").Value = "0.0₃0793"
$ws.Range("E39").Value = "  -5.53%  "
$ws.Range("D40").Value = "3.364.84"
$ws.Range("E40").Value = "  +10.00%  "
$ws.Range("E41").Value = "  -6.17%  "
$ws.Range("D42").Value = "'3.28"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("E44").Value = "  -4.98%  "
$ws.Range("D45").Value = "'2.54"
$ws.Range("E45").Value = "  -10.02%  "
$ws.Range("D46").Value = "'3.28"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'0.0414"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'138.26"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.46"
$ws.Range("E51").Value = "  -3.68%  "
